# Generate Report for Handoff
# Update Priority to "ht" and refresh the Latest Handoff Datetime for the
# files that were just generated for handoff (rows 4-7: 00732b48, 2c2af40c,
# 5a920252, 63f6a715) on both the zh-cn and de-de localization sheets.

$wb = $excel.ActiveWorkbook

$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("E4:E7").Value = "ht"
$zhcn.Range("H4:H7").Value = "2016-08-26 10:29:30"

$dede = $wb.Worksheets.Item("de-de")
$dede.Range("E4:E7").Value = "ht"
$dede.Range("H4:H7").Value = "2016-08-26 10:29:36"

# The Overview sheet's "Latest HO Xliff Generate Date" column mirrors the
# de-de handoff timestamp for the same file, so refresh it too.
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G4:G7").Value = "2016-08-26 10:29:36"
